# khl_stats_1369_ext.xlsx — publish new matches + refresh derived stats
# (chore(runtime): publish files + archive (2025-11-12 17:26:40))

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Matches_SOG: append the two newly-played matches (rows 463 & 464)
# ---------------------------------------------------------------------------
$matches = $wb.Worksheets.Item("Matches_SOG")

# Column A (uid) holds numeric-looking values but is stored as text in the
# source data, so force text formatting before writing to avoid Excel's
# automatic number coercion.
$matches.Range("A463").NumberFormat = "@"
$matches.Range("A464").NumberFormat = "@"

$matches.Range("A463").Value = "897760"
$matches.Range("B463").Value = "2025-11-11T12:15:00"
$matches.Range("C463").Value = "Амур"
$matches.Range("D463").Value = "Трактор"
$matches.Range("E463").Value = 27
$matches.Range("F463").Value = 36
$matches.Range("G463").Value = "khl_text"

$matches.Range("A464").Value = "897761"
$matches.Range("B464").Value = "2025-11-11T12:30:00"
$matches.Range("C464").Value = "Адмирал"
$matches.Range("D464").Value = "Автомобилист"
$matches.Range("E464").Value = 33
$matches.Range("F464").Value = 21
$matches.Range("G464").Value = "khl_text"

# ---------------------------------------------------------------------------
# 2) Shots_HA: bump as_of_utc on every team row, refresh the four teams that
#    actually played (Автомобилист / Адмирал / Амур / Трактор).
# ---------------------------------------------------------------------------
$shotsHA = $wb.Worksheets.Item("Shots_HA")
for ($r = 2; $r -le 23; $r++) {
    $shotsHA.Range("D$r").Value = "2025-11-11T12:30:00Z"
}

# row 3 = Автомобилист
$shotsHA.Range("F3").Value = 27
$shotsHA.Range("K3").Value = 752
$shotsHA.Range("L3").Value = 831
$shotsHA.Range("M3").Value = 27.9
$shotsHA.Range("N3").Value = 30.8

# row 4 = Адмирал
$shotsHA.Range("E4").Value = 18
$shotsHA.Range("G4").Value = 686
$shotsHA.Range("H4").Value = 483
$shotsHA.Range("I4").Value = 38.1
$shotsHA.Range("J4").Value = 26.8

# row 6 = Амур
$shotsHA.Range("E6").Value = 21
$shotsHA.Range("G6").Value = 641
$shotsHA.Range("H6").Value = 746
$shotsHA.Range("I6").Value = 30.5

# row 21 = Трактор
$shotsHA.Range("F21").Value = 26
$shotsHA.Range("K21").Value = 888
$shotsHA.Range("L21").Value = 847
$shotsHA.Range("M21").Value = 34.2
$shotsHA.Range("N21").Value = 32.6

# ---------------------------------------------------------------------------
# 3) Shots_Summary: same idea — bump as_of_utc everywhere, refresh totals for
#    the four teams that played.
# ---------------------------------------------------------------------------
$shotsSummary = $wb.Worksheets.Item("Shots_Summary")
for ($r = 2; $r -le 23; $r++) {
    $shotsSummary.Range("D$r").Value = "2025-11-11T12:30:00Z"
}

# row 3 = Автомобилист
$shotsSummary.Range("E3").Value = 45
$shotsSummary.Range("F3").Value = 1273
$shotsSummary.Range("G3").Value = 1389
$shotsSummary.Range("H3").Value = 28.3
$shotsSummary.Range("I3").Value = 30.9

# row 4 = Адмирал
$shotsSummary.Range("E4").Value = 38
$shotsSummary.Range("F4").Value = 1324
$shotsSummary.Range("G4").Value = 1044
$shotsSummary.Range("H4").Value = 34.8
$shotsSummary.Range("I4").Value = 27.5

# row 6 = Амур
$shotsSummary.Range("E6").Value = 42
$shotsSummary.Range("F6").Value = 1233
$shotsSummary.Range("G6").Value = 1515

# row 21 = Трактор
$shotsSummary.Range("E21").Value = 44
$shotsSummary.Range("F21").Value = 1487
$shotsSummary.Range("G21").Value = 1388
$shotsSummary.Range("H21").Value = 33.8
$shotsSummary.Range("I21").Value = 31.5

# ---------------------------------------------------------------------------
# 4) Meta_ext: bump as_of_utc + build_version.
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Meta_ext")
$meta.Range("B2").Value = "2025-11-11T12:30:00Z"
$meta.Range("D2").Value = 60
